# Actualización automática 2025-09-11 08:50:10
# Insert a new advisor row ("PUEBLA GONZALEZ MARIO DANIEL") before the
# existing "TAMAYO VILLACIS EDWIN XAVIER" row (row 19) on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets. All following rows shift
# down by one; the new row carries the same advisor name in column A and
# zeroes across the numeric columns. The trailing summary row's
# "X de 22" counters become "X de 23" to reflect the new total headcount.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R, data rows 2-23, summary row 24)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row above row 19 (shifts old rows 19-24 down to 20-25).
$ws1.Rows.Item(19).Insert()

$ws1.Cells.Item(19, 1).Value = "RIOS CARRION ANGEL BENIGNO"
$ws1.Cells.Item(19, 2).Value = "PUEBLA GONZALEZ MARIO DANIEL"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(19, $col).Value = 0
}

# Update the "X de 22" -> "X de 23" summary labels on the (new) row 25.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(25, $col)
    $text = $cell.Value()
    $text = $text.Replace("de 22", "de 23")
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G, data rows 2-23, summary row 24)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a blank row above row 19 (shifts old rows 19-24 down to 20-25).
$ws2.Rows.Item(19).Insert()

$ws2.Cells.Item(19, 1).Value = "RIOS CARRION ANGEL BENIGNO"
$ws2.Cells.Item(19, 2).Value = "PUEBLA GONZALEZ MARIO DANIEL"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(19, $col).Value = 0
}
